# Applies the odds updates described in the commit diff to Sheet1.
# Three match rows (3, 6, 8) receive updated odds values in various columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Alianza vs Fortaleza) ---
$ws.Range("G3").Value  = 2.1
$ws.Range("I3").Value  = 3.8
$ws.Range("AB3").Value = 9
$ws.Range("AI3").Value = 17
$ws.Range("AL3").Value = 17
$ws.Range("AO3").Value = 34
$ws.Range("AP3").Value = 41

# --- Row 6 (Racing Montevideo vs River Plate) ---
$ws.Range("G6").Value  = 1.83
$ws.Range("H6").Value  = 3
$ws.Range("I6").Value  = 4.5
$ws.Range("J6").Value  = 2.5
$ws.Range("K6").Value  = 2.1
$ws.Range("Y6").Value  = 1.8
$ws.Range("Z6").Value  = 1.91
$ws.Range("AC6").Value = 9
$ws.Range("AH6").Value = 6
$ws.Range("AM6").Value = 15

# --- Row 8 (Plaza Colonia vs Boston River) ---
$ws.Range("G8").Value  = 2
$ws.Range("H8").Value  = 3.3
$ws.Range("I8").Value  = 3.5
$ws.Range("J8").Value  = 2.88
$ws.Range("L8").Value  = 4.75
$ws.Range("M8").Value  = 1.1
$ws.Range("N8").Value  = 7
$ws.Range("O8").Value  = 1.5
$ws.Range("P8").Value  = 2.5
$ws.Range("Q8").Value  = 1.98
$ws.Range("R8").Value  = 1.88
$ws.Range("S8").Value  = 2.5
$ws.Range("T8").Value  = 1.5
$ws.Range("Y8").Value  = 2.25
$ws.Range("Z8").Value  = 1.57
$ws.Range("AA8").Value = 5.5
$ws.Range("AB8").Value = 8.5
$ws.Range("AC8").Value = 10
$ws.Range("AD8").Value = 17
$ws.Range("AE8").Value = 21
$ws.Range("AG8").Value = 7
$ws.Range("AH8").Value = 7
$ws.Range("AK8").Value = 7.5
$ws.Range("AL8").Value = 17
$ws.Range("AM8").Value = 13
$ws.Range("AN8").Value = 41
$ws.Range("AO8").Value = 41
$ws.Range("AP8").Value = 51
